$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete first match row (Mexican Liga MX: Santos Laguna vs FC Juarez).
# Deleting it shifts the remaining match rows up by one, matching the diff.
$ws.Rows.Item(2).Delete()

# Update row 2 with refreshed match data
$ws.Cells.Item(2, 1).Value = "Colombian Primera A"
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 2).Value = "2026-01-18"
$ws.Cells.Item(2, 3).Value = "20:20:00"
$ws.Cells.Item(2, 4).Value = "Junior FC Barranquilla"
$ws.Cells.Item(2, 5).Value = "Tolima"
$ws.Cells.Item(2, 6).Value = 44
$ws.Cells.Item(2, 7).Value = 1000
$ws.Cells.Item(2, 8).Value = 1.07
$ws.Cells.Item(2, 9).Value = 1.08
$ws.Cells.Item(2, 10).Value = 15.5
$ws.Cells.Item(2, 11).Value = 1000
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = 0
$ws.Cells.Item(2, 14).Value = 0
$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(2, 16).Value = 1.66
$ws.Cells.Item(2, 17).Value = 2.46
$ws.Cells.Item(2, 18).Value = 1.07
$ws.Cells.Item(2, 19).Value = 7.6
$ws.Cells.Item(2, 20).Value = 1.1
$ws.Cells.Item(2, 21).Value = 1.02
$ws.Cells.Item(2, 22).Value = 7.6
$ws.Cells.Item(2, 23).Value = 1.01
$ws.Cells.Item(2, 24).Value = 1000
$ws.Cells.Item(2, 25).Value = 1000
$ws.Cells.Item(2, 26).Value = 1.69
$ws.Cells.Item(2, 27).Value = 1000
$ws.Cells.Item(2, 28).Value = 1000
$ws.Cells.Item(2, 29).Value = 1000
$ws.Cells.Item(2, 30).Value = 1000
$ws.Cells.Item(2, 31).Value = 1000
$ws.Cells.Item(2, 32).Value = 1000
$ws.Cells.Item(2, 33).Value = 1000
$ws.Cells.Item(2, 34).Value = 1000
$ws.Cells.Item(2, 35).Value = 1000
$ws.Cells.Item(2, 36).Value = 1000
$ws.Cells.Item(2, 37).Value = 1000
$ws.Cells.Item(2, 38).Value = 1000
$ws.Cells.Item(2, 39).Value = 1000
$ws.Cells.Item(2, 40).Value = 1000
$ws.Cells.Item(2, 41).Value = 980

# Update row 3 with refreshed match data
$ws.Cells.Item(3, 1).Value = "Mexican Liga MX"
$ws.Cells.Item(3, 2).NumberFormat = "@"
$ws.Cells.Item(3, 2).Value = "2026-01-18"
$ws.Cells.Item(3, 3).Value = "22:06:00"
$ws.Cells.Item(3, 4).Value = "Pachuca"
$ws.Cells.Item(3, 5).Value = "CF America"
$ws.Cells.Item(3, 6).Value = 3.2
$ws.Cells.Item(3, 7).Value = 3.25
$ws.Cells.Item(3, 8).Value = 2.62
$ws.Cells.Item(3, 9).Value = 2.64
$ws.Cells.Item(3, 10).Value = 3.25
$ws.Cells.Item(3, 11).Value = 3.3
$ws.Cells.Item(3, 12).Value = 1.62
$ws.Cells.Item(3, 13).Value = 1.12
$ws.Cells.Item(3, 14).Value = 2.92
$ws.Cells.Item(3, 15).Value = 1.5
$ws.Cells.Item(3, 16).Value = 1.59
$ws.Cells.Item(3, 17).Value = 2.6
$ws.Cells.Item(3, 18).Value = 1.22
$ws.Cells.Item(3, 19).Value = 5.2
$ws.Cells.Item(3, 20).Value = 2
$ws.Cells.Item(3, 21).Value = 1.87
$ws.Cells.Item(3, 22).Value = 1.61
$ws.Cells.Item(3, 23).Value = 1.45
$ws.Cells.Item(3, 24).Value = 9.199999999999999
$ws.Cells.Item(3, 25).Value = 8.199999999999999
$ws.Cells.Item(3, 26).Value = 16
$ws.Cells.Item(3, 27).Value = 46
$ws.Cells.Item(3, 28).Value = 9.6
$ws.Cells.Item(3, 29).Value = 7
$ws.Cells.Item(3, 30).Value = 13.5
$ws.Cells.Item(3, 31).Value = 38
$ws.Cells.Item(3, 32).Value = 21
$ws.Cells.Item(3, 33).Value = 15
$ws.Cells.Item(3, 34).Value = 23
$ws.Cells.Item(3, 35).Value = 70
$ws.Cells.Item(3, 36).Value = 75
$ws.Cells.Item(3, 37).Value = 48
$ws.Cells.Item(3, 38).Value = 70
$ws.Cells.Item(3, 39).Value = 210
$ws.Cells.Item(3, 40).Value = 65
$ws.Cells.Item(3, 41).Value = 40

# Update row 4 with refreshed match data
$ws.Cells.Item(4, 1).Value = "Colombian Primera A"
$ws.Cells.Item(4, 2).NumberFormat = "@"
$ws.Cells.Item(4, 2).Value = "2026-01-18"
$ws.Cells.Item(4, 3).Value = "22:30:00"
$ws.Cells.Item(4, 4).Value = "Santa Fe"
$ws.Cells.Item(4, 5).Value = "Aguilas Doradas"
$ws.Cells.Item(4, 6).Value = 1.99
$ws.Cells.Item(4, 7).Value = 2.04
$ws.Cells.Item(4, 8).Value = 4.5
$ws.Cells.Item(4, 9).Value = 4.8
$ws.Cells.Item(4, 10).Value = 3.4
$ws.Cells.Item(4, 11).Value = 3.55
$ws.Cells.Item(4, 12).Value = 1.49
$ws.Cells.Item(4, 13).Value = 1.1
$ws.Cells.Item(4, 14).Value = 3.2
$ws.Cells.Item(4, 15).Value = 1.43
$ws.Cells.Item(4, 16).Value = 1.75
$ws.Cells.Item(4, 17).Value = 2.26
$ws.Cells.Item(4, 18).Value = 1.29
$ws.Cells.Item(4, 19).Value = 4.4
$ws.Cells.Item(4, 20).Value = 2
$ws.Cells.Item(4, 21).Value = 1.92
$ws.Cells.Item(4, 22).Value = 1.27
$ws.Cells.Item(4, 23).Value = 1.96
$ws.Cells.Item(4, 24).Value = 11
$ws.Cells.Item(4, 25).Value = 13.5
$ws.Cells.Item(4, 26).Value = 32
$ws.Cells.Item(4, 27).Value = 110
$ws.Cells.Item(4, 28).Value = 8
$ws.Cells.Item(4, 29).Value = 7.8
$ws.Cells.Item(4, 30).Value = 18.5
$ws.Cells.Item(4, 31).Value = 70
$ws.Cells.Item(4, 32).Value = 11
$ws.Cells.Item(4, 33).Value = 10.5
$ws.Cells.Item(4, 34).Value = 21
$ws.Cells.Item(4, 35).Value = 85
$ws.Cells.Item(4, 36).Value = 23
$ws.Cells.Item(4, 37).Value = 23
$ws.Cells.Item(4, 38).Value = 46
$ws.Cells.Item(4, 39).Value = 150
$ws.Cells.Item(4, 40).Value = 19
$ws.Cells.Item(4, 41).Value = 90
